$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values in column A (rows 2-6)
$ws.Range("A2").Value = 11111
$ws.Range("A3").Value = 22222
$ws.Range("A4").Value = 33333
$ws.Range("A5").Value = 44444
$ws.Range("A6").Value = 55555

# Update the selected cell/range to A6
$ws.Range("A6").Select()
